$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append the two new daily-data rows (129 and 130) ------------------
# Seed formatting by copying the last existing data row, then overwrite values.
$ws.Range("A128:G128").Copy()
$ws.Range("A129:G129").PasteSpecial(-4122)
$ws.Range("A128:G128").Copy()
$ws.Range("A130:G130").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

$ws.Range("A129").Value = 44050
$ws.Range("B129").Value = 12
$ws.Range("C129").Value = 1720
$ws.Range("D129").Value = 80
$ws.Range("E129").Value = 1515
$ws.Range("F129").Value = 1073
$ws.Range("G129").Value = 4

$ws.Range("A130").Value = 44051
$ws.Range("B130").Value = 11
$ws.Range("C130").Value = 1731
$ws.Range("D130").Value = 80
$ws.Range("E130").Value = 1529
$ws.Range("F130").Value = 1080
$ws.Range("G130").Value = 2

# --- 2. Point the first chart's data series at the new, shifted range -----
$co1 = $ws.ChartObjects(1)
$chart1 = $co1.Chart
$ser1 = $chart1.SeriesCollection(1)
$ser1.Formula = "=SERIES(,,Planilha1!`$B`$4:`$B`$130,1)"

# --- 3. Reposition/resize both chart objects on the sheet ------------------
$co1.Top = 2165.812440944882
$co1.Left = 772.6260356176181
$co1.Width = 802.8124606299214
$co1.Height = 231.5625196850392

$co2 = $ws.ChartObjects(2)
$co2.Top = 2414.062362204724
$co2.Left = 772.6257993971457
$co2.Width = 817.1877
$co2.Height = 125.81267716535467

# --- 4. Restore the active selection left by the editing session ----------
$ws.Range("F171").Select()
